$wb = $excel.ActiveWorkbook

$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Range("G2").Value = 33
$sprint1.Range("G3").Value = 30
$sprint1.Range("G4").Value = 12
$sprint1.Range("G5").Value = 13
$sprint1.Range("G6").Value = 26
$sprint1.Range("G7").Value = 28

$burndown = $wb.Worksheets.Item("Burndown")
$burndown.Range("A2").Value = 40953
$burndown.Range("B2").Value = 42
$burndown.Range("C2").Value = 0
$burndown.Range("D2").Value = 0
$burndown.Range("E2").Value = 0
$burndown.Range("F2").Value = 0

$burndown.Range("A3").Value = 40966
$burndown.Range("B3").Formula = "=42-6"
$burndown.Range("C3").Formula = "=B2-B3"
$burndown.Range("D3").Formula = "=SUM(Sprint1!G2:G7)"
$burndown.Range("E3").Formula = "=MIN(Sprint1!G2:G7)"
$burndown.Range("F3").Formula = "=(D3-D2)/E3*60"

$co = $burndown.ChartObjects().Item(1)
$chart = $co.Chart
try {
    $chart.SetSourceData($burndown.Range("A1:B7"))
    Write-Output "setsourcedata ok"
} catch {
    Write-Output "setsourcedata failed: $_"
}
